$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-02-26", "sleep", $true, $false),
    @("2025-02-26", "activity", $false, $false),
    @("2025-02-26", "weekly_activity", $false, $false)
)

$startRow = 77
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    # Force the date column to be stored as plain text, matching the
    # existing rows (which store dates as literal strings, not real
    # Excel date values), then reset the style so no extra formatting
    # is left behind on the cell.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
